$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the obsolete "Unnamed: 0" column (old column B). This
#    shifts venue_id/venue/venue_code/.../geometry one column to the
#    left (old C:M -> new B:L) for every existing row, including the
#    header.
# ------------------------------------------------------------------
$ws.Range("B1").EntireColumn.Delete()

# ------------------------------------------------------------------
# 2. Rename venue "Ward Community Center" -> "Ward Town Hall"
#    (row 31, venue name now lives in column C after the shift).
# ------------------------------------------------------------------
$ws.Range("C31").Value = "Ward Town Hall"

# ------------------------------------------------------------------
# 3. Fill in the previously-blank "Meadowlark" venue row (row 35,
#    venue index 33) - the Erie venue moved to Meadowlark.
# ------------------------------------------------------------------
$ws.Range("B35").Value = 33
$ws.Range("C35").Value = "Meadowlark"
$ws.Range("D35").Value = "MWLK_P8"
$ws.Range("E35").Value = "2300 Meadow Sweet Ln, Erie, CO 80516"
$ws.Range("F35").Value = "BVSD"
$ws.Range("G35").Value = "https://ml8.bvsd.org/"
$ws.Range("H35").Value = "https://maps.app.goo.gl/EA1rf9bvvJL3sp8aA"
$ws.Range("I35").Value = "Meadowlark School"
$ws.Range("J35").Value = 40.034318124098
$ws.Range("K35").Value = -105.083013187148
$ws.Range("L35").Value = "POINT (-105.083013187148 40.034318124098)"

# ------------------------------------------------------------------
# 4. Row 37 no longer carries the stray lat/lon note that used to
#    sit in column H (it was shifted there from the old I37 when the
#    column was deleted above). That note now belongs to a brand new
#    row 38, with coordinates refreshed for the Meadowlark venue and
#    at higher precision ("Judi's Forecast" note row).
# ------------------------------------------------------------------
$ws.Range("H37").ClearContents()

$ws.Range("A37").Copy($ws.Range("A38"))
$ws.Range("A38").Value = 36
$ws.Range("I38").Value = "40.03431812409801, -105.08301318714825"
$ws.Range("L38").Value = "POINT (nan nan)"
